# Refresh the cryptos price/volume snapshot (GitHub Actions style data update).
# Numeric-looking price strings (e.g. "586.28") are prefixed with a leading
# apostrophe so Excel stores them as text, matching the source data's
# inlineStr cell type instead of auto-coercing to a number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.298.96'
$ws.Range('E2').Value = '  +1.15%  '
$ws.Range('D3').Value = '3.502.77'
$ws.Range('E3').Value = '  +0.80%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '''586.28'
$ws.Range('E5').Value = '  +0.97%  '
$ws.Range('D6').Value = '''134.29'
$ws.Range('E6').Value = '  +3.66%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('E8').Value = '  +1.24%  '
$ws.Range('E9').Value = '  +1.97%  '
$ws.Range('E10').Value = '  +2.56%  '
$ws.Range('D11').Value = '''0.386'
$ws.Range('E11').Value = '  +2.94%  '
$ws.Range('D12').Value = '4.099.77'
$ws.Range('E12').Value = '  +1.03%  '
$ws.Range('E13').Value = '  +4.32%  '
$ws.Range('E14').Value = '  +1.50%  '
$ws.Range('D15').Value = '3.502.35'
$ws.Range('E15').Value = '  +0.14%  '
$ws.Range('D16').Value = '''26.04'
$ws.Range('E16').Value = '  -3.68%  '
$ws.Range('D17').Value = '64.313.05'
$ws.Range('E17').Value = '  +1.07%  '
$ws.Range('D18').Value = '''9.90'
$ws.Range('E18').Value = '  +0.99%  '
$ws.Range('E19').Value = '  +3.31%  '
$ws.Range('D20').Value = '''13.65'
$ws.Range('E20').Value = '  -2.29%  '
$ws.Range('D21').Value = '''393.69'
$ws.Range('E21').Value = '  +4.09%  '
$ws.Range('E22').Value = '  +0.45%  '
$ws.Range('D23').Value = '3.644.16'
$ws.Range('E23').Value = '  +0.75%  '
$ws.Range('D24').Value = '''74.29'
$ws.Range('E24').Value = '  +1.82%  '
$ws.Range('E25').Value = '  +0.05%  '
$ws.Range('E26').Value = '  +3.00%  '
$ws.Range('E27').Value = '  +0.03%  '
$ws.Range('D28').Value = '''7.41'
$ws.Range('E28').Value = '  +0.25%  '
$ws.Range('B29').Value = 'InternetComputer(DFINITY)'
$ws.Range('C29').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D29').Value = '''8.30'
$ws.Range('E29').Value = '  +1.85%  '
$ws.Range('B30').Value = 'Fetch.AI'
$ws.Range('C30').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D30').Value = '''1.50'
$ws.Range('E30').Value = '  -3.60%  '
$ws.Range('E31').Value = '  +1.67%  '
$ws.Range('D32').Value = '3.523.57'
$ws.Range('E32').Value = '  +1.06%  '
$ws.Range('E33').Value = '  +5.09%  '
$ws.Range('D35').Value = '''23.45'
$ws.Range('E35').Value = '  +0.90%  '
$ws.Range('D36').Value = '''5.16'
$ws.Range('E36').Value = '  -1.19%  '
$ws.Range('D37').Value = '''1.56'
$ws.Range('E37').Value = '  +1.61%  '
$ws.Range('E38').Value = '  +0.67%  '
$ws.Range('D39').Value = '''164.02'
$ws.Range('E39').Value = '  +2.76%  '
$ws.Range('E40').Value = '  +0.03%  '
$ws.Range('E41').Value = '  +0.12%  '
$ws.Range('E42').Value = '  +0.06%  '
$ws.Range('D43').Value = '''25.31'
$ws.Range('E43').Value = '  -2.03%  '
$ws.Range('D44').Value = '''4.41'
$ws.Range('E44').Value = '  +2.23%  '
$ws.Range('D45').Value = '''1.66'
$ws.Range('E45').Value = '  +3.94%  '
$ws.Range('E46').Value = '  -1.51%  '
$ws.Range('D47').Value = '2.464.26'
$ws.Range('E47').Value = '  +2.24%  '
$ws.Range('E48').Value = '  -0.05%  '
$ws.Range('D49').Value = '''0.899'
$ws.Range('E49').Value = '  +1.98%  '
$ws.Range('E50').Value = '  +0.00%  '
$ws.Range('E51').Value = '  +0.36%  '
